{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"Consolida\u00e7\u00e3o das Leis do Trabalho...\" reference paragraph, which stays intact.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"Consolida\u00e7\u00e3o das Leis do Trabalho\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find the anchor paragraph.\");\n}\n\n// Immediately following the anchor are exactly three paragraphs to remove:\n//   1) a blank \"Normal\" paragraph\n//   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3) \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n// Verify the expected texts line up before removing anything, then delete from the\n// bottom up so earlier indices stay valid.\nconst expectedSnippets = [\"\", \"Ver no Jupiter\", \"Powered by Jekyll and Github pages\"];\nconst candidates = items.slice(anchorIndex + 1, anchorIndex + 1 + expectedSnippets.length);\n\nif (candidates.length !== expectedSnippets.length) {\n  throw new Error(\"Not enough trailing paragraphs to match the expected pattern.\");\n}\n\nfor (let i = 0; i < expectedSnippets.length; i++) {\n  const snippet = expectedSnippets[i];\n  const text = candidates[i].text;\n  const matches = snippet === \"\" ? text === \"\" : text.indexOf(snippet) !== -1;\n  if (!matches) {\n    throw new Error(\"Paragraph \" + i + \" did not match expected content: \" + text);\n  }\n}\n\nfor (let i = candidates.length - 1; i >= 0; i--) {\n  candidates[i].delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"Consolida\u00e7\u00e3o das Leis do Trabalho...\" reference paragraph, which stays intact.\n$anchorIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text\n    if ($text -like \"*Consolida*Leis do Trabalho*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find the anchor paragraph.\"\n}\n\n# Immediately following the anchor are exactly three paragraphs to remove:\n#   1) a blank \"Normal\" paragraph\n#   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   3) \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n$targets = @($anchorIndex + 1, $anchorIndex + 2, $anchorIndex + 3)\n\n$expected = @(\n    \"\",\n    \"Ver no Jupiter\",\n    \"Powered by Jekyll and Github pages\"\n)\n\nfor ($j = 0; $j -lt $targets.Length; $j++) {\n    $idx = $targets[$j]\n    $text = $d.Paragraphs.Item($idx).Range.Text.TrimEnd([char]13, [char]7)\n    $snippet = $expected[$j]\n    if ($snippet -eq \"\") {\n        if ($text -ne \"\") {\n            throw \"Paragraph $idx did not match expected blank content: $text\"\n        }\n    } elseif ($text -notlike \"*$snippet*\") {\n        throw \"Paragraph $idx did not match expected content: $text\"\n    }\n}\n\n# Delete from the bottom up so earlier indices stay valid.\nfor ($j = $targets.Length - 1; $j -ge 0; $j--) {\n    $idx = $targets[$j]\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
